$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16-79 hold "Periodo Mora" (E) / "Valor Mora" (F) records.
# The account-statement periods are reordered to newest-first (2110 .. 1607)
# while each period keeps the mora value it originally had.
$periods = @("2110","2109","2108","2107","2106","2105","2104","2103","2102","2101","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001","1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901","1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801","1712","1711","1710","1709","1708","1707","1706","1705","1704","1703","1702","1701","1612","1611","1610","1609","1608","1607")
$values = @(26041,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $values[$i]
}
